# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# handback (ede20ffa... and f8a17299... files) has completed for both the
# zh-cn and de-de locales:
#   - Status text changes from "In Translation" to
#     "Handed back: in sync with en-US" everywhere it appears.
#   - The per-locale "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns are populated.
#   - Column widths are widened to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

$mdEdeUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e17b722f9c88d2d0a6fbfc7312213f7b7649cba1/e2e/ede20ffa-7bdb-40ca-a178-139789da6e0c.md"
$mdF8aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e17b722f9c88d2d0a6fbfc7312213f7b7649cba1/e2e/f8a17299-f943-4972-86d1-abf4af00a71d.md"
$mdEdeName = "ede20ffa-7bdb-40ca-a178-139789da6e0c.md"
$mdF8aName = "f8a17299-f943-4972-86d1-abf4af00a71d.md"

# ---------------------------------------------------------------------
# 1. Overview sheet: Status columns (zh-cn = E, de-de = F) for both rows
# ---------------------------------------------------------------------
$ws1.Range("E2").Value = $handedBack
$ws1.Range("F2").Value = $handedBack
$ws1.Range("E3").Value = $handedBack
$ws1.Range("F3").Value = $handedBack

$ws1.Columns.Item(5).ColumnWidth = 29.15
$ws1.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# 2. zh-cn sheet (table "zh_cn")
#    Columns: C=Status, I=Latest Target File, J=Latest Handback File,
#             K=Latest Handback DateTime
# ---------------------------------------------------------------------
$ws2.Range("C2").Value = $handedBack
$ws2.Range("C3").Value = $handedBack

$ws2.Range("I2").Value = $mdEdeName
$ws2.Range("J2").Value = "ede20ffa-7bdb-40ca-a178-139789da6e0c.c7d21694720db5eb8a1371a1af5923b83a91d8f4.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-10-14 08:39:52"

$ws2.Range("I3").Value = $mdF8aName
$ws2.Range("J3").Value = "f8a17299-f943-4972-86d1-abf4af00a71d.44d0c2ed93eb6d581e55139912b056af1342d5f7.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-10-14 08:39:52"

$ws2.Hyperlinks.Add($ws2.Range("I2"), $mdEdeUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdEdeName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), $mdF8aUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdF8aName) | Out-Null

$ws2.Range("I2").Font.Underline = $true
$ws2.Range("I2").Font.Color = 15631332
$ws2.Range("I3").Font.Underline = $true
$ws2.Range("I3").Font.Color = 15631332

$ws2.Columns.Item(3).ColumnWidth = 29.15
$ws2.Columns.Item(9).ColumnWidth = 39.15
$ws2.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# 3. de-de sheet (table "de_de")
#    Columns: C=Status, I=Latest Target File, J=Latest Handback File,
#             K=Latest Handback DateTime
# ---------------------------------------------------------------------
$ws3.Range("C2").Value = $handedBack
$ws3.Range("C3").Value = $handedBack

$ws3.Range("I2").Value = $mdEdeName
$ws3.Range("J2").Value = "ede20ffa-7bdb-40ca-a178-139789da6e0c.c7d21694720db5eb8a1371a1af5923b83a91d8f4.de-de.xlf"
$ws3.Range("K2").Value = "2016-10-14 08:40:10"

$ws3.Range("I3").Value = $mdF8aName
$ws3.Range("J3").Value = "f8a17299-f943-4972-86d1-abf4af00a71d.44d0c2ed93eb6d581e55139912b056af1342d5f7.de-de.xlf"
$ws3.Range("K3").Value = "2016-10-14 08:40:10"

$ws3.Hyperlinks.Add($ws3.Range("I2"), $mdEdeUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdEdeName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), $mdF8aUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdF8aName) | Out-Null

$ws3.Range("I2").Font.Underline = $true
$ws3.Range("I2").Font.Color = 15631332
$ws3.Range("I3").Font.Underline = $true
$ws3.Range("I3").Font.Color = 15631332

$ws3.Columns.Item(3).ColumnWidth = 29.15
$ws3.Columns.Item(9).ColumnWidth = 39.15
$ws3.Columns.Item(10).ColumnWidth = 39.15

Write-Output "Handback report generated."
